# edit.ps1 - apply the tracked changes to CW2_CLASS_DIAGRAM_TEMPLATE_TEAM 9.docx
$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $ok = $range.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "NOT FOUND: $find"
    }
}

# 1) "requirements listed above..." paragraph: remove the proofErr-wrapped "not" run split
#    (text content unchanged - "and not other " stays the same; nothing to replace)

# 2) "runs the whole application" -> "runs the application"
Replace-Text "runs the whole application" "runs the application"

# 3) private lists description rewrite
Replace-Text "private lists of admins, students, quizzes, and categories, and uses public methods" "private lists which are admins, students, quizzes and categories, it also uses public methods"

# 4) remove comma before "and its list of"
Replace-Text "number of questions, and its list of" "number of questions and its list of"

# 5) remove comma before "and manage quizzes"
Replace-Text "display results, and manage quizzes" "display results and manage quizzes"

# 6) add "also" before "provides methods"
Replace-Text "correct answer and provides methods" "correct answer and also provides methods"

# 7) Results -> results (lower-case) + restructure "The User class information" -> "While the User class stores information"
Replace-Text "while the Results class records a user" "while the results class records a user"
Replace-Text "wrong answers. The User class information like username, password, email, and role, with methods" "wrong answers. While the User class stores information like username, password, email and role and with methods"

# 8) role-specific behaviours. -> unique classes.  (also merges in the trailing "." run
#    and adds a trailing space so the old standalone period run disappears)
Replace-Text "role-specific behaviours." "unique classes. "

# 9) The stray paragraph that used to contain only "." becomes fully empty
#    (its content merges into the end of the previous paragraph as "unique classes. ")
$pDot = $d.Paragraphs.Item(75)
if ($pDot.Range.Text.TrimEnd([char]13) -eq ".") {
    $rDot = $pDot.Range
    $rDot.MoveEnd(1, -1) | Out-Null
    $rDot.Text = ""
}

Write-Output "Done with text replacements"
